{"js": "// Replace \"Currently, I am working with\" with \"Currently, I was working with\"\n// (matches the author's edit: \"am\" -> \"was\" inside the bio paragraph).\nconst body = context.document.body;\nconst results = body.search(\"Currently, I am working with\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found\");\n}\n\nconst found = results.items[0];\nfound.insertText(\"Currently, I was working with\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Currently, I am working\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"Currently, I was working\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Execute([ref]\"Currently, I am working\", [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, \"Currently, I was working\", 2) | Out-Null\n"}
